# Update metricas_retencao_anual: refreshed num_customers (and cohort_size for
# the last row) values after reloading data for ADD/BIBI cohorts, with the
# dependent retention_rate recalculated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# row, C (num_customers), D (cohort_size, only set when it changes)
$updates = @(
    @{ Row = 22; C = 24;  D = $null },
    @{ Row = 27; C = 36;  D = $null },
    @{ Row = 31; C = 36;  D = $null },
    @{ Row = 34; C = 58;  D = $null },
    @{ Row = 36; C = 99;  D = $null },
    @{ Row = 37; C = 584; D = 584 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value2 = $u.C
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value2 = $u.D
    }
    $numCustomers = $ws.Cells.Item($r, 3).Value2
    $cohortSize = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value2 = $numCustomers / $cohortSize
}
